$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.786.02'
$ws.Range("E2").Value = '  -1.76%  '

# Row 3
$ws.Range("D3").Value = '2.353.07'
$ws.Range("E3").Value = '  -1.29%  '

# Row 4
$ws.Range("E4").Value = '  -0.24%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.28'
$ws.Range("E5").Value = '  +1.51%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.61'
$ws.Range("E6").Value = '  -4.13%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.640'
$ws.Range("E7").Value = '  +0.10%  '

# Row 8
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.626'
$ws.Range("E9").Value = '  -1.95%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.17'
$ws.Range("E10").Value = '  -5.48%  '

# Row 11
$ws.Range("E11").Value = '  -1.69%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.50'
$ws.Range("E12").Value = '  -2.69%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.00'
$ws.Range("E13").Value = '  -4.62%  '

# Row 14
$ws.Range("E14").Value = '  +0.25%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.18'
$ws.Range("E15").Value = '  -3.01%  '

# Row 16
$ws.Range("D16").Value = '2.708.46'
$ws.Range("E16").Value = '  -1.25%  '

# Row 17
$ws.Range("D17").Value = '2.349.08'
$ws.Range("E17").Value = '  -4.57%  '

# Row 18
$ws.Range("D18").Value = '42.723.75'
$ws.Range("E18").Value = '  -1.97%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.89'
$ws.Range("E19").Value = '  +8.44%  '

# Row 20
$ws.Range("E20").Value = '  -2.40%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.84'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.67'
$ws.Range("E22").Value = '  +4.76%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '264.12'
$ws.Range("E23").Value = '  +0.89%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.32'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.01'
$ws.Range("E25").Value = '  +8.38%  '

# Row 26
$ws.Range("E26").Value = '  -0.03%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.50'
$ws.Range("E27").Value = '  -4.69%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.95'
$ws.Range("E28").Value = '  +0.03%  '

# Row 29
$ws.Range("E29").Value = '  -1.34%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.90'
$ws.Range("E30").Value = '  +0.40%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.12'
$ws.Range("E31").Value = '  -3.07%  '

# Row 32
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0899'
$ws.Range("E32").Value = '  -2.81%  '

# Row 33
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.41'
$ws.Range("E33").Value = '  -8.90%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.15'
$ws.Range("E34").Value = '  +2.90%  '

# Row 35
$ws.Range("E35").Value = '  +0.99%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.112'
$ws.Range("E36").Value = '  +6.61%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.56'
$ws.Range("E37").Value = '  -8.11%  '

# Row 38
$ws.Range("E38").Value = '  -4.20%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.80'
$ws.Range("E39").Value = '  -7.12%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.71'
$ws.Range("E40").Value = '  -4.69%  '

# Row 41
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.237'
$ws.Range("E41").Value = '  +1.58%  '

# Row 42
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.48'
$ws.Range("E42").Value = '  -1.96%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.26'
$ws.Range("E43").Value = '  -2.40%  '

# Row 44
$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '94.89'
$ws.Range("E44").Value = '  +25.95%  '

# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '121.16'
$ws.Range("E45").Value = '  +7.21%  '

# Row 46
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  -0.21%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.90'
$ws.Range("E47").Value = '  -5.41%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.54'
$ws.Range("E48").Value = '  -1.73%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.10'
$ws.Range("E49").Value = '  -2.62%  '

# Row 50
$ws.Range("E50").Value = '  -3.84%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.100'
$ws.Range("E51").Value = '  +0.08%  '
